$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "24.868.04"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +1.69%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.708.91"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.75%  "

$ws.Range("E4").Value = "  +0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "312.95"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.99%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3756"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.51%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "49.60"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.09%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3471"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.53%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.220"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.44%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07582"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.29%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.02%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.42"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.06%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.351"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +3.35%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.095"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.17%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.711.52"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.03%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001136"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.42%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.06739"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.29%  "

$ws.Range("E19").Value = "  -0.03%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "85.04"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +4.82%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.40"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.79%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.423"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +5.18%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "13.24"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +10.17%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "24.895.69"
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.461"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.805"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.88%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.56"
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "151.09"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "133.24"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +5.13%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.901.67"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.11%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.250"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +27.87%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.921"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +9.06%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "4.236"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +4.67%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "13.96"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +10.87%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.08870"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +5.02%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.766"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +3.74%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.675"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +5.78%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "9.371"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +5.14%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.06695"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +2.91%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.02424"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +4.18%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.2251"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +6.46%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.284"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.93%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6495"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.09%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.9990"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("E45").Value = "  +6.16%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6190"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.94%  "

$ws.Range("E47").Value = "  +2.11%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.145"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +5.70%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "130.90"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.81%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.07320"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.59%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "80.58"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.05%  "

